# Refresh cryptocurrency price/volume data scraped from coinranking.com
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '27.140.16'
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').Value = '1.900.92'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.36%  '
$ws.Range('E6').Value = '  +0.32%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5233'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.43%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3806'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.86%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07287'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.60%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.38'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.07%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9037'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.69%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08217'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.52%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '95.41'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.74%  '
$ws.Range('D14').Value = '1.842.77'
$ws.Range('E14').Value = '  -3.10%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.356'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.004'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.34%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008676'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.58%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.68'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.04%  '
$ws.Range('E19').Value = '  +0.30%  '
$ws.Range('D20').Value = '27.175.42'
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('B22').Value = 'Cosmos'
$ws.Range('C22').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.80'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.10%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.446'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.41%  '
$ws.Range('B24').Value = 'Monero'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '149.86'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.39%  '
$ws.Range('B25').Value = 'LidoDAOToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.324'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.87%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.29'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.05%  '
$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.740'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.94%  '
$ws.Range('B28').Value = 'BitcoinCash'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '115.66'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.82%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.820'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.89%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.902'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.35%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09221'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.45%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05046'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.00%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7938'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.82%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.222'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.03%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.965'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.39%  '
$ws.Range('B36').Value = 'MXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.360'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.36%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.644'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.96%  '
$ws.Range('B38').Value = 'TheSandbox'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5730'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.31%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01992'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.65%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.082'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.87%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '9.097'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.69%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.612'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.79%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '116.41'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.38%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1517'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.40%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4890'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.18%  '
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.004'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.39%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.14'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.55%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.637'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.51%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '38.52'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.91%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '63.96'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.59%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05956'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.52%  '
